$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 952.64703
$ws.Range("I19").Value = 573.75
$ws.Range("J19").Value = 1289.4445
$ws.Range("K19").Value = 573.75
$ws.Range("L19").Value = 1289.4445
$ws.Range("M19").Value = -398.75
$ws.Range("N19").Value = -1639.4445

$ws.Range("H88").Value = 1600.6
$ws.Range("I88").Value = 1267.6666
$ws.Range("J88").Value = 2100
$ws.Range("K88").Value = 1267.6666
$ws.Range("L88").Value = 2100
$ws.Range("M88").Value = -861.6666
$ws.Range("N88").Value = -2912

$ws.Range("H91").Value = 1600.6
$ws.Range("I91").Value = 1267.6666
$ws.Range("J91").Value = 2100
$ws.Range("K91").Value = 1267.6666
$ws.Range("L91").Value = 2100
$ws.Range("M91").Value = 136.3334
$ws.Range("N91").Value = -4908

$ws.Range("H106").Value = 3645.4
$ws.Range("I106").Value = 2188.5
$ws.Range("K106").Value = 2188.5
$ws.Range("M106").Value = -1557.5

$ws.Range("H112").Value = 4802.1
$ws.Range("I112").Value = 14352
$ws.Range("J112").Value = 1618.8
$ws.Range("K112").Value = 43056
$ws.Range("L112").Value = 4856.4
$ws.Range("M112").Value = -41948
$ws.Range("N112").Value = -7072.4

$ws.Range("H115").Value = 2891.4
$ws.Range("I115").Value = 2891.4
$ws.Range("K115").Value = 8674.200000000001
$ws.Range("M115").Value = -7107.200000000001

$ws.Range("H129").Value = 801.1277
$ws.Range("I129").Value = 345.33334
$ws.Range("J129").Value = 909.0789
$ws.Range("K129").Value = 1036.00002
$ws.Range("L129").Value = 2727.2367
$ws.Range("M129").Value = 3963.99998
$ws.Range("N129").Value = -12727.2367

$ws.Range("H137").Value = 3256.524
$ws.Range("I137").Value = 3594.353
$ws.Range("J137").Value = 1820.75
$ws.Range("K137").Value = 10783.059
$ws.Range("L137").Value = 5462.25
$ws.Range("M137").Value = -8233.059000000001
$ws.Range("N137").Value = -10562.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25890.705
$ws.Range("I32").Value = 25890.705
$ws.Range("K32").Value = 25890.705
$ws.Range("M32").Value = -25603.705

$ws.Range("H45").Value = 1811.7059
$ws.Range("I45").Value = 1878.5714
$ws.Range("J45").Value = 1499.6666
$ws.Range("K45").Value = 1878.5714
$ws.Range("L45").Value = 1499.6666
$ws.Range("M45").Value = -1501.5714
$ws.Range("N45").Value = -2253.6666

$ws.Range("H63").Value = 3247.4
$ws.Range("I63").Value = 2497.7144
$ws.Range("J63").Value = 4996.6665
$ws.Range("K63").Value = 2497.7144
$ws.Range("L63").Value = 4996.6665
$ws.Range("M63").Value = -1811.7144
$ws.Range("N63").Value = -6368.6665

$ws.Range("H66").Value = 3247.4
$ws.Range("I66").Value = 2497.7144
$ws.Range("J66").Value = 4996.6665
$ws.Range("K66").Value = 12488.572
$ws.Range("L66").Value = 24983.3325
$ws.Range("M66").Value = -9056.572
$ws.Range("N66").Value = -31847.3325

$ws.Range("H74").Value = 2398.0322
$ws.Range("I74").Value = 2092.3809
$ws.Range("J74").Value = 3039.9
$ws.Range("K74").Value = 2092.3809
$ws.Range("L74").Value = 3039.9
$ws.Range("M74").Value = -1218.3809
$ws.Range("N74").Value = -4787.9

$ws.Range("H77").Value = 2398.0322
$ws.Range("I77").Value = 2092.3809
$ws.Range("J77").Value = 3039.9
$ws.Range("K77").Value = 10461.9045
$ws.Range("L77").Value = 15199.5
$ws.Range("M77").Value = -6093.904500000001
$ws.Range("N77").Value = -23935.5

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("N101").ClearContents()

$ws.Range("H132").Value = 7548.864
$ws.Range("I132").Value = 9467.462
$ws.Range("J132").Value = 4777.5557
$ws.Range("K132").Value = 28402.386
$ws.Range("L132").Value = 14332.6671
$ws.Range("M132").Value = -25872.386
$ws.Range("N132").Value = -19392.6671

$ws.Range("H134").Value = 53832.332
$ws.Range("J134").Value = 53832.332
$ws.Range("L134").Value = 53832.332
$ws.Range("N134").Value = -63972.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 311.4
$ws.Range("I80").Value = 147.5
$ws.Range("J80").Value = 352.375
$ws.Range("K80").Value = 147.5
$ws.Range("L80").Value = 352.375
$ws.Range("M80").Value = 850.5
$ws.Range("N80").Value = -2348.375

$ws.Range("H83").Value = 311.4
$ws.Range("I83").Value = 147.5
$ws.Range("J83").Value = 352.375
$ws.Range("K83").Value = 737.5
$ws.Range("L83").Value = 1761.875
$ws.Range("M83").Value = 4254.5
$ws.Range("N83").Value = -11745.875

$ws.Range("H105").Value = 3577.2258
$ws.Range("I105").Value = 3562.1904
$ws.Range("J105").Value = 3608.8
$ws.Range("K105").Value = 3562.1904
$ws.Range("L105").Value = 3608.8
$ws.Range("M105").Value = -1815.1904
$ws.Range("N105").Value = -7102.8

$ws.Range("H107").Value = 2255.1853
$ws.Range("I107").Value = 1866.6666
$ws.Range("K107").Value = 1866.6666
$ws.Range("M107").Value = 53.33339999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 12041.1
$ws.Range("I16").Value = 17735.166
$ws.Range("K16").Value = 17735.166
$ws.Range("M16").Value = -17448.166

$ws.Range("H31").Value = 9093.166999999999
$ws.Range("I31").Value = 8048.357
$ws.Range("K31").Value = 8048.357
$ws.Range("M31").Value = -7753.357

$ws.Range("H34").Value = 9093.166999999999
$ws.Range("I34").Value = 8048.357
$ws.Range("K34").Value = 8048.357
$ws.Range("M34").Value = -7846.357

$ws.Range("H58").Value = 2935083
$ws.Range("I58").Value = 6994941
$ws.Range("J58").Value = 2963.2222
$ws.Range("K58").Value = 6994941
$ws.Range("L58").Value = 2963.2222
$ws.Range("M58").Value = -6994738
$ws.Range("N58").Value = -3369.2222

$ws.Range("H86").Value = 4500
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 4500
$ws.Range("N86").Value = -6746
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 4500
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 22500
$ws.Range("N89").Value = -33732
$ws.Range("M89").ClearContents()

$ws.Range("H99").Value = 2427.2727
$ws.Range("I99").Value = 1883.3334
$ws.Range("J99").Value = 3080
$ws.Range("K99").Value = 1883.3334
$ws.Range("L99").Value = 3080
$ws.Range("M99").Value = -385.3334
$ws.Range("N99").Value = -6076

$ws.Range("H110").Value = 79350
$ws.Range("J110").Value = 79350
$ws.Range("L110").Value = 79350
$ws.Range("N110").Value = -87530

$ws.Range("H113").Value = 12041.1
$ws.Range("I113").Value = 17735.166
$ws.Range("K113").Value = 17735.166
$ws.Range("M113").Value = -15565.166

$ws.Range("H122").Value = 3589.5652
$ws.Range("I122").Value = 3288.05
$ws.Range("J122").Value = 5599.6665
$ws.Range("K122").Value = 9864.150000000001
$ws.Range("L122").Value = 16798.9995
$ws.Range("M122").Value = -7414.150000000001
$ws.Range("N122").Value = -21698.9995

$ws.Range("H126").Value = 2427.2727
$ws.Range("I126").Value = 1883.3334
$ws.Range("J126").Value = 3080
$ws.Range("K126").Value = 5650.0002
$ws.Range("L126").Value = 9240
$ws.Range("M126").Value = -3180.0002
$ws.Range("N126").Value = -14180

$ws.Range("H132").Value = 3340.524
$ws.Range("I132").Value = 2710.7693
$ws.Range("J132").Value = 4363.875
$ws.Range("K132").Value = 8132.3079
$ws.Range("L132").Value = 13091.625
$ws.Range("M132").Value = -5602.3079
$ws.Range("N132").Value = -18151.625

$ws.Range("H136").Value = 2935083
$ws.Range("I136").Value = 6994941
$ws.Range("J136").Value = 2963.2222
$ws.Range("K136").Value = 20984823
$ws.Range("L136").Value = 8889.6666
$ws.Range("M136").Value = -20982273
$ws.Range("N136").Value = -13989.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1203.9231
$ws.Range("J122").Value = 1633.7778
$ws.Range("L122").Value = 14704.0002
$ws.Range("N122").Value = -19604.0002

$ws.Range("H131").Value = 881.56665
$ws.Range("I131").Value = 967.1177
$ws.Range("J131").Value = 847.7442
$ws.Range("K131").Value = 2901.3531
$ws.Range("L131").Value = 2543.2326
$ws.Range("M131").Value = 2138.6469
$ws.Range("N131").Value = -12623.2326

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3940.111
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 4410.1665
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 4410.1665
$ws.Range("M40").Value = -2864
$ws.Range("N40").Value = -4682.1665

$ws.Range("H105").Value = 40615
$ws.Range("J105").Value = 40615
$ws.Range("L105").Value = 40615
$ws.Range("N105").Value = -47603

$ws.Range("H136").Value = 5633.9707
$ws.Range("I136").Value = 3225.3572
$ws.Range("J136").Value = 7320
$ws.Range("K136").Value = 9676.071599999999
$ws.Range("L136").Value = 21960
$ws.Range("M136").Value = -7126.071599999999
$ws.Range("N136").Value = -27060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 47100.5
$ws.Range("J103").Value = 47100.5
$ws.Range("L103").Value = 47100.5
$ws.Range("N103").Value = -49444.5

$ws.Range("H113").Value = 727.1
$ws.Range("J113").Value = 997.375
$ws.Range("L113").Value = 2992.125
$ws.Range("N113").Value = -7332.125

$ws.Range("H132").Value = 4884.154
$ws.Range("I132").Value = 5685.2856
$ws.Range("J132").Value = 3949.5
$ws.Range("K132").Value = 17055.8568
$ws.Range("L132").Value = 11848.5
$ws.Range("M132").Value = -14525.8568
$ws.Range("N132").Value = -16908.5
